# Updates cryptos list values (Price / Volume(1h) columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the literal string into the cell without Excel's automatic
    # number/date inference (values like "291.94" or "2.45" would
    # otherwise be parsed as numbers), then restore the cell's original
    # (default) style so no formatting changes are introduced.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "39.909.34"
$ws.Range("E2").Value = "  +0.40%  "
Set-TextValue "D3" "2.225.25"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "291.94"
$ws.Range("E5").Value = "  -0.18%  "
Set-TextValue "D6" "87.29"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.82%  "
Set-TextValue "D10" "30.52"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -0.58%  "
Set-TextValue "D12" "49.99"
$ws.Range("E12").Value = "  +5.57%  "
$ws.Range("E13").Value = "  +2.65%  "
Set-TextValue "D14" "6.45"
$ws.Range("E14").Value = "  +1.62%  "
Set-TextValue "D15" "2.568.12"
$ws.Range("E15").Value = "  +0.69%  "
Set-TextValue "D16" "13.84"
$ws.Range("E16").Value = "  -1.43%  "
Set-TextValue "D17" "2.233.55"
$ws.Range("E17").Value = "  +0.99%  "
Set-TextValue "D18" "0.733"
$ws.Range("E18").Value = "  +0.59%  "
Set-TextValue "D19" "39.833.33"
$ws.Range("E19").Value = "  +0.37%  "
Set-TextValue "D20" "0.0₃0887"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("E22").Value = "  -0.80%  "
Set-TextValue "D23" "65.85"
$ws.Range("E23").Value = "  +0.20%  "
Set-TextValue "D24" "237.21"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("E25").Value = "  -0.05%  "
Set-TextValue "D26" "2.46"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  -0.56%  "
Set-TextValue "D28" "23.07"
$ws.Range("E28").Value = "  +1.48%  "
Set-TextValue "D29" "9.25"
$ws.Range("E29").Value = "  -0.25%  "
Set-TextValue "D31" "156.88"
$ws.Range("E31").Value = "  +3.01%  "
Set-TextValue "D32" "31.97"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  +0.03%  "
Set-TextValue "D34" "4.97"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("E35").Value = "  +6.98%  "
Set-TextValue "D36" "0.0715"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  +3.20%  "
Set-TextValue "D40" "0.0993"
$ws.Range("E40").Value = "  +0.48%  "
Set-TextValue "D41" "15.33"
$ws.Range("E41").Value = "  -4.03%  "
Set-TextValue "D42" "2.112.04"
$ws.Range("E42").Value = "  +1.58%  "
Set-TextValue "D43" "3.73"
$ws.Range("E43").Value = "  -1.75%  "
Set-TextValue "D44" "18.17"
$ws.Range("E44").Value = "  +2.22%  "
Set-TextValue "D45" "0.0272"
$ws.Range("E45").Value = "  +1.17%  "
Set-TextValue "D46" "9.94"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("E47").Value = "  -7.77%  "
Set-TextValue "D48" "2.73"
$ws.Range("E48").Value = "  +4.68%  "
Set-TextValue "D49" "2.441.44"
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  +2.42%  "
